# feat: split tipo recebimento from setor
#
# Updates the "entidade" sheet CNPJ value, bumps the FTP folder paths
# from the "...1" suite to the "...2" suite, and moves the active
# sheet/selection back to "entidade".

$wb = $excel.ActiveWorkbook

$entidade = $wb.Worksheets.Item("entidade")
$ftp      = $wb.Worksheets.Item("ftp")
$padrao   = $wb.Worksheets.Item("padrao_integracao")

# entidade: new CNPJ + matching numeric inscricao_estadual
$entidade.Range("C2").Value = "12.345.678/0001-05"
$entidade.Range("E2").Value = 123456700145

# ftp: bump the "1" suffixed paths to "2"
$ftp.Range("A2").Value = "/home/ftpsynapcomp/Embu/3M/importacao2"
$ftp.Range("B2").Value = "/home/ftpsynapcomp/Embu/3M/bkp_importacao2"
$ftp.Range("C2").Value = "/home/ftpsynapcomp/Embu/3M/exportacao2"
$ftp.Range("D2").Value = "/home/ftpsynapcomp/Embu/3M/bkp_exportacao2"
$ftp.Range("E2").Value = "/home/ftpsynapcomp/Embu/3M/erro2"

# Move the selection / active sheet back to entidade
[void]$entidade.Range("D17").Select()
[void]$entidade.Activate()

# Update the lingering selections on the other two touched sheets
[void]$ftp.Range("A19").Select()
[void]$padrao.Range("G30").Select()

[void]$entidade.Activate()
